$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capacity")

# --- Update capacities / electricity demand values for 2030 ---
$ws.Range("F101").Value = 15
$ws.Range("H102").Value = 900
$ws.Range("F103").Value = 20
$ws.Range("H104").Value = 750
$ws.Range("F105").Value = 10

# --- Re-style rows 116-118 (FR00 Electrolysis / Hydrogen storage dimensioner /
#     Hydrogen processor) so columns A-C match the highlighted format already
#     used in columns D-H ---
$ws.Range("D116").Copy()
$ws.Range("A116:C116").PasteSpecial(-4122)
$ws.Range("D117").Copy()
$ws.Range("A117:C117").PasteSpecial(-4122)
$ws.Range("D118").Copy()
$ws.Range("A118:C118").PasteSpecial(-4122)

# --- Add new row 173: FI00 / Hydrogen processor / Distributed Energy / 2030 / 300 ---
$ws.Range("A173").Value = "FI00"
$ws.Range("C173").Value = "Hydrogen processor"
$ws.Range("D173").Value = "Distributed Energy"
$ws.Range("E173").Value = 2030
$ws.Range("H173").Value = 300
$ws.Range("A116:H116").Copy()
$ws.Range("A173:H173").PasteSpecial(-4122)

# --- Apply additional AutoFilter on the Year column (E, colId 4) to 2030,
#     on top of the existing Generator_ID filter. This hides every row whose
#     Year isn't 2030 (rows that were previously visible because they matched
#     the Generator_ID filter). ---
$ws.Range("A1:J173").AutoFilter(5, @("2030"), 7)

# --- Restore the cursor position recorded in the saved file ---
$ws.Range("F176").Select()
